$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BoM")

# --- Header metadata block (rows 6-9, column C) -------------------------
# These cells used to hold literal sample values copied from the BomDoc
# report ("Virtus CC", the doc filename, the version, the date). They are
# switched to the generic "Field = <Name>" placeholders that the BoM
# report generator substitutes at export time.
$ws.Range("C6").Value = "Field = Projeto"
$ws.Range("C7").Value = "Field = DataSourceFileName"
$ws.Range("C8").Value = "Field = Rev"
$ws.Range("C9").Value = "Field = ReportDate"

# --- Table header row (row 12) ------------------------------------------
# The column headers now carry the "Column=<Name>" placeholder syntax used
# by the report generator (and the previous, broken
# "#Column Name Error:' Manufacturer Part Number" label is fixed).
$ws.Range("C12").Value = "Column=Quantity"
$ws.Range("D12").Value = "Column=Designator"
$ws.Range("E12").Value = "Column=Manufacturer Part Number"
$ws.Range("F12").Value = "Column=Description"
$ws.Range("G12").Value = "Column=Supplier 1"
$ws.Range("H12").Value = "Column=Supplier Part Number 1"
$ws.Range("I12").Value = "Column=Supplier Stock 1"
$ws.Range("J12").Value = "Column=Supplier Unit Price 1"
$ws.Range("K12").Value = "Column=Supplier Subtotal 1"
